# Update the cached regression-table figures that changed when the
# underlying cleaning files / internal-validity analysis for the
# "number of pawns" balance table were re-run (decomposition_main_te).
#
# The sheet's cells are formulas pulling cached results from an external
# workbook ([1]decomposition_main_te!...); we refresh the handful of
# cells whose figures moved by typing the new value straight into the
# cell, exactly as the source table was updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (source row 5): main coefficients
$ws.Range("B6").Value = "-204.0***"
$ws.Range("E6").Value = "-1.27"
$ws.Range("F6").Value = "-78.8**"
$ws.Range("G6").Value = "-0.066***"

# Row 7 (source row 6): standard errors
$ws.Range("B7").Value = "(48.1)"
$ws.Range("E7").Value = "(3.10)"

# Row 8 (source row 8): second set of coefficients
$ws.Range("B8").Value = "-38.9"
$ws.Range("E8").Value = "-0.93"
$ws.Range("F8").Value = "-15.4"
$ws.Range("I8").Value = "-0.0086"

# Row 9 (source row 9): standard errors
$ws.Range("E9").Value = "(3.02)"
$ws.Range("F9").Value = "(33.1)"

# Row 13 (source row 13): control mean
$ws.Range("B13").Value = "942.4"
$ws.Range("E13").Value = "5.96"
$ws.Range("F13").Value = "396.5"
$ws.Range("G13").Value = "0.44"
